# Edit the "8.0gb memory is not enough (t3a.large) ..." bullet so it reads
# "16.0gb memory seems not enough (t3a.xlarge) ..." while reproducing the
# run layout of the target revision (six runs, each with an explicit empty
# <w:rPr/>, split at the points where the wording changed).

$d = $word.ActiveDocument

# Locate the paragraph that still holds the original wording.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "8.0gb memory is not enough*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the paragraph starting with '8.0gb memory is not enough'."
}

$start = $target.Range.Start

# --- Step 1: fix up the wording --------------------------------------------
# Apply edits back-to-front so earlier character offsets stay valid while the
# paragraph length changes underneath us.

# "(t3a.large)" -> "(t3a.xlarge)" : insert an "x" right before "large"
$rX = $d.Range($start + 32, $start + 32)
$rX.Text = "x"

# "is" -> "seems"
$rIs = $d.Range($start + 13, $start + 15)
$rIs.Text = "seems"

# "8" -> "16"
$r8 = $d.Range($start, $start + 1)
$r8.Text = "16"

# --- Step 2: re-split the run into the six fragments seen in the target ----
# Toggling a character-formatting property on and back off forces a run
# break at the range boundaries without altering the visible formatting
# (it reverts to the original/inherited value): the engine only re-merges
# adjacent runs when run *text* actually changes, not when formatting is
# merely (re)applied.

$bounds = @(0, 2, 14, 19, 36, 37, 74)
for ($i = 0; $i -lt $bounds.Length - 1; $i++) {
    $seg = $d.Range($start + $bounds[$i], $start + $bounds[$i + 1])
    $seg.Bold = 1
    $seg.Bold = 0
}
